$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Tir_235_50R24")
$ws.Rows.Item(8).Insert()
$ws.Cells.Item(8, 1).Value = "roadFile"
$ws.Cells.Item(8, 8).Value = "which('TNO_FlatRoad.rdf')"
